$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for the specified rows
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F9").Value = -4
$ws.Range("F11").Value = 4
$ws.Range("F16").Value = 0
$ws.Range("F18").Value = 2
$ws.Range("F23").Value = 5
$ws.Range("F25").Value = -4
